# The underlying source data lost one record (Caso 3659, OT 797489950,
# "ALBERTI 59") that had been included in a previous export. Removing that
# row shifts every subsequent record up by one position on both the
# "General" sheet (full list) and the "PEBCOM" sheet (filtered view for
# that provider), and the used range shrinks by one row on each.

$wb = $excel.ActiveWorkbook

# "General" sheet: the record is on spreadsheet row 6 (5th data row).
$wsGeneral = $wb.Worksheets.Item("General")
$wsGeneral.Rows.Item(6).Delete()

# "PEBCOM" sheet: the same record, filtered, is on spreadsheet row 4.
$wsPebcom = $wb.Worksheets.Item("PEBCOM")
$wsPebcom.Rows.Item(4).Delete()
